$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Data table: Row, C (N Doc Trabajador), D (Nombre Trabajador), E (Periodo Mora), F (Valor Mora), G (Salario Basico)
$data = @(
    @(16, '45452851', 'YADIRA TRIBIÑO CORTINA', '2208', 38133, 1300000),
    @(17, '45452851', 'YADIRA TRIBIÑO CORTINA', '2207', 52000, 1300000),
    @(18, '45452851', 'YADIRA TRIBIÑO CORTINA', '2206', 52000, 1300000),
    @(19, '45452851', 'YADIRA TRIBIÑO CORTINA', '2205', 52000, 1300000),
    @(20, '45452851', 'YADIRA TRIBIÑO CORTINA', '2204', 52000, 1300000),
    @(21, '45452851', 'YADIRA TRIBIÑO CORTINA', '2203', 52000, 1300000),
    @(22, '45506395', 'BEATRIZ ELENA SALAZAR HERRERA', '2208', 52454, 1788211),
    @(23, '45506395', 'BEATRIZ ELENA SALAZAR HERRERA', '2207', 71528, 1788211),
    @(24, '45506395', 'BEATRIZ ELENA SALAZAR HERRERA', '2206', 71528, 1788211),
    @(25, '45506395', 'BEATRIZ ELENA SALAZAR HERRERA', '2205', 71528, 1788211),
    @(26, '45506395', 'BEATRIZ ELENA SALAZAR HERRERA', '2204', 71528, 1788211),
    @(27, '45506395', 'BEATRIZ ELENA SALAZAR HERRERA', '2203', 71528, 1788211),
    @(28, '45506395', 'BEATRIZ ELENA SALAZAR HERRERA', '2202', 71528, 1788211),
    @(29, '32909140', 'ERICA PATRICIA SOTO MENDOZA', '2202', 30284, 1000000),
    @(30, '45689665', 'JOHANNA DELFINA CHAVEZ MORA', '2208', 29333, 0),
    @(31, '45689665', 'JOHANNA DELFINA CHAVEZ MORA', '2207', 40000, 0),
    @(32, '45689665', 'JOHANNA DELFINA CHAVEZ MORA', '2206', 40000, 0),
    @(33, '45689665', 'JOHANNA DELFINA CHAVEZ MORA', '2205', 40000, 0),
    @(34, '45689665', 'JOHANNA DELFINA CHAVEZ MORA', '2204', 40000, 0),
    @(35, '45689665', 'JOHANNA DELFINA CHAVEZ MORA', '2203', 40000, 0),
    @(36, '45689665', 'JOHANNA DELFINA CHAVEZ MORA', '2202', 80000, 0),
    @(37, '73105496', 'JAVIER JESUS COLL JIMENEZ', '2208', 29333, 1000000),
    @(38, '73105496', 'JAVIER JESUS COLL JIMENEZ', '2207', 40000, 1000000),
    @(39, '73105496', 'JAVIER JESUS COLL JIMENEZ', '2206', 40000, 1000000),
    @(40, '73105496', 'JAVIER JESUS COLL JIMENEZ', '2205', 40000, 1000000),
    @(41, '73105496', 'JAVIER JESUS COLL JIMENEZ', '2204', 40000, 1000000),
    @(42, '73105496', 'JAVIER JESUS COLL JIMENEZ', '2203', 40000, 1000000),
    @(43, '73105496', 'JAVIER JESUS COLL JIMENEZ', '2202', 40000, 1000000)
)

foreach ($row in $data) {
    $r    = $row[0]
    $docC = $row[1]
    $nomD = $row[2]
    $perE = $row[3]
    $valF = $row[4]
    $salG = $row[5]

    $ws.Cells.Item($r, 3).Value = $docC
    $ws.Cells.Item($r, 4).Value = $nomD
    $ws.Cells.Item($r, 5).Value = $perE
    $ws.Cells.Item($r, 6).Value = $valF
    $ws.Cells.Item($r, 7).Value = $salG
}
